# This sheet is a weekly price log for "Poroto verde" (Macroferia Regional de
# Talca). The edit inserts one new weekly record as row 126, pushing the
# existing rows 126-154 down to 127-155 (so the former last row, 154, becomes
# the new last row, 155). Columns A, B, C, E, F, G, H, I, N, Q, R are constant
# across every record in this sheet, so the new row reuses those values; only
# the date (D), volume (J), prices (K/L/M), origin (O) and $/Kg (P) are new.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 126; Excel shifts rows 126:154 down to 127:155 and
# carries their formatting (incl. the date style on column D) along with them.
$ws.Rows("126:126").Insert()

# Populate the newly inserted row 126 with the new weekly record.
$ws.Cells.Item(126, 1).Value2 = 5
$ws.Cells.Item(126, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(126, 3).Value2 = "Maule"
$ws.Cells.Item(126, 4).Value2 = 44641
$ws.Cells.Item(126, 5).Value2 = 7
$ws.Cells.Item(126, 6).Value2 = 100112031
$ws.Cells.Item(126, 7).Value2 = "Poroto verde"
$ws.Cells.Item(126, 8).Value2 = "Sin especificar"
$ws.Cells.Item(126, 9).Value2 = "Primera"
$ws.Cells.Item(126, 10).Value2 = 200
$ws.Cells.Item(126, 11).Value2 = 23000
$ws.Cells.Item(126, 12).Value2 = 23000
$ws.Cells.Item(126, 13).Value2 = 23000
$ws.Cells.Item(126, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(126, 15).Value2 = "Región del Maule"
$ws.Cells.Item(126, 16).Value2 = 920
$ws.Cells.Item(126, 17).Value2 = 25
$ws.Cells.Item(126, 18).Value2 = "Hortaliza"
